$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 140.209918
$ws.Range("H2").Value = 420.629754
$ws.Range("I2").Value = 0.2353423205412711
$ws.Range("J2").Value = 0.2363790708159033
$ws.Range("M2").Value = 8.753652666666667
$ws.Range("N2").Value = 26.260958
$ws.Range("O2").Value = 0.2520793895170918
$ws.Range("P2").Value = 0.2812533920401061
$ws.Range("Q2").Value = 1227.348922593815
$ws.Range("R2").Value = 11046.14030334433
$ws.Range("S2").Value = 0.05932494848957934
$ws.Range("T2").Value = 0.06648241547426124
$ws.Range("G3").Value = 140.209918
$ws.Range("H3").Value = 420.629754
$ws.Range("I3").Value = 0.2353423205412711
$ws.Range("J3").Value = 0.2363790708159033
$ws.Range("O3").Value = 0.4367350033185243
$ws.Range("P3").Value = 0.4872798261741802
$ws.Range("Q3").Value = 2126.418335147755
$ws.Range("R3").Value = 19137.7650163298
$ws.Range("S3").Value = 0.1027822291425812
$ws.Range("T3").Value = 0.1151827525383876
$ws.Range("G4").Value = 140.209918
$ws.Range("H4").Value = 420.629754
$ws.Range("I4").Value = 0.2353423205412711
$ws.Range("J4").Value = 0.2363790708159033
$ws.Range("M4").Value = 10.806162
$ws.Range("N4").Value = 21.612324
$ws.Range("O4").Value = 0.311185607164384
$ws.Range("P4").Value = 0.2314667817857137
$ws.Range("Q4").Value = 1515.131087914716
$ws.Range("R4").Value = 9090.786527488297
$ws.Range("S4").Value = 0.07323514290911051
$ws.Range("T4").Value = 0.05471390280325446
$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "ECs"
$ws.Range("G5").Value = 446.6683856666667
$ws.Range("H5").Value = 1340.005157
$ws.Range("I5").Value = 0.7497327998952026
$ws.Range("J5").Value = 0.753035587444864
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 8.753652666666667
$ws.Range("N5").Value = 26.260958
$ws.Range("O5").Value = 0.2520793895170918
$ws.Range("P5").Value = 0.2812533920401061
$ws.Range("Q5").Value = 3909.979905306712
$ws.Range("R5").Value = 35189.81914776041
$ws.Range("S5").Value = 0.1889921864985226
$ws.Range("T5").Value = 0.2117938132957819
$ws.Range("D6").Value = "FAPs"
$ws.Range("I6").Value = 0.7497327998952026
$ws.Range("J6").Value = 0.753035587444864
$ws.Range("M6").Value = 15.16596233333333
$ws.Range("N6").Value = 45.49788700000001
$ws.Range("O6").Value = 0.4367350033185243
$ws.Range("P6").Value = 0.4872798261741802
$ws.Range("Q6").Value = 6774.155912511474
$ws.Range("R6").Value = 60967.40321260327
$ws.Range("S6").Value = 0.3274345568502378
$ws.Range("T6").Value = 0.366939050153105
$ws.Range("D7").Value = "MuSCs"
$ws.Range("I7").Value = 0.7497327998952026
$ws.Range("J7").Value = 0.753035587444864
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 10.806162
$ws.Range("N7").Value = 21.612324
$ws.Range("O7").Value = 0.311185607164384
$ws.Range("P7").Value = 0.2314667817857137
$ws.Range("Q7").Value = 4826.770935792478
$ws.Range("R7").Value = 28960.62561475487
$ws.Range("S7").Value = 0.2333060565464422
$ws.Range("T7").Value = 0.1743027239959771
$ws.Range("A8").Value = "Inflammatory-Mac"
$ws.Range("D8").Value = "ECs"
$ws.Range("G8").Value = 0.7424606666666667
$ws.Range("H8").Value = 2.227382
$ws.Range("I8").Value = 0.001246220086969543
$ws.Range("J8").Value = 0.001251710043108525
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 8.753652666666667
$ws.Range("N8").Value = 26.260958
$ws.Range("O8").Value = 0.2520793895170918
$ws.Range("P8").Value = 0.2812533920401061
$ws.Range("Q8").Value = 6.499242794661778
$ws.Range("R8").Value = 58.49318515195601
$ws.Range("S8").Value = 0.0003141463987272194
$ws.Range("T8").Value = 0.0003520476954749401
$ws.Range("A9").Value = "Inflammatory-Mac"
$ws.Range("D9").Value = "FAPs"
$ws.Range("G9").Value = 0.7424606666666667
$ws.Range("H9").Value = 2.227382
$ws.Range("I9").Value = 0.001246220086969543
$ws.Range("J9").Value = 0.001251710043108525
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 15.16596233333333
$ws.Range("N9").Value = 45.49788700000001
$ws.Range("O9").Value = 0.4367350033185243
$ws.Range("P9").Value = 0.4872798261741802
$ws.Range("Q9").Value = 11.26013050464822
$ws.Range("R9").Value = 101.341174541834
$ws.Range("S9").Value = 0.000544267933818255
$ws.Range("T9").Value = 0.0006099330522263977
$ws.Range("D10").Value = "MuSCs"
$ws.Range("G10").Value = 0.7424606666666667
$ws.Range("H10").Value = 2.227382
$ws.Range("I10").Value = 0.001246220086969543
$ws.Range("J10").Value = 0.001251710043108525
$ws.Range("K10").Value = 2
$ws.Range("M10").Value = 10.806162
$ws.Range("N10").Value = 21.612324
$ws.Range("O10").Value = 0.311185607164384
$ws.Range("P10").Value = 0.2314667817857137
$ws.Range("Q10").Value = 8.023150242628001
$ws.Range("R10").Value = 48.138901455768
$ws.Range("S10").Value = 0.0003878057544240686
$ws.Range("T10").Value = 0.0002897292954071873
$ws.Range("A11").Value = "MuSCs"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 2
$ws.Range("G11").Value = 7.839080000000001
$ws.Range("H11").Value = 15.67816
$ws.Range("I11").Value = 0.01315789428040795
$ws.Range("J11").Value = 0.008810572380248361
$ws.Range("M11").Value = 8.753652666666667
$ws.Range("N11").Value = 26.260958
$ws.Range("O11").Value = 0.2520793895170918
$ws.Range("P11").Value = 0.2812533920401061
$ws.Range("Q11").Value = 68.62058354621335
$ws.Range("R11").Value = 411.7235012772801
$ws.Range("S11").Value = 0.00331683395753567
$ws.Range("T11").Value = 0.002478003367759723
$ws.Range("A12").Value = "MuSCs"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 2
$ws.Range("G12").Value = 7.839080000000001
$ws.Range("H12").Value = 15.67816
$ws.Range("I12").Value = 0.01315789428040795
$ws.Range("J12").Value = 0.008810572380248361
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 15.16596233333333
$ws.Range("N12").Value = 45.49788700000001
$ws.Range("O12").Value = 0.4367350033185243
$ws.Range("P12").Value = 0.4872798261741802
$ws.Range("Q12").Value = 118.8871920079867
$ws.Range("R12").Value = 713.3231520479202
$ws.Range("S12").Value = 0.005746513002218757
$ws.Range("T12").Value = 0.004293214177942454
$ws.Range("A13").Value = "MuSCs"
$ws.Range("D13").Value = "MuSCs"
$ws.Range("E13").Value = 2
$ws.Range("G13").Value = 7.839080000000001
$ws.Range("H13").Value = 15.67816
$ws.Range("I13").Value = 0.01315789428040795
$ws.Range("J13").Value = 0.008810572380248361
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 10.806162
$ws.Range("N13").Value = 21.612324
$ws.Range("O13").Value = 0.311185607164384
$ws.Range("P13").Value = 0.2314667817857137
$ws.Range("Q13").Value = 84.71036841096002
$ws.Range("R13").Value = 338.8414736438401
$ws.Range("S13").Value = 0.004094547320653523
$ws.Range("T13").Value = 0.002039354834546184
$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.3102563333333333
$ws.Range("H14").Value = 0.930769
$ws.Range("I14").Value = 0.0005207651961489113
$ws.Range("J14").Value = 0.0005230593158758034
$ws.Range("M14").Value = 8.753652666666667
$ws.Range("N14").Value = 26.260958
$ws.Range("O14").Value = 0.2520793895170918
$ws.Range("P14").Value = 0.2812533920401061
$ws.Range("Q14").Value = 2.715876179633555
$ws.Range("R14").Value = 24.442885616702
$ws.Range("S14").Value = 0.0001312741727269661
$ws.Range("T14").Value = 0.000147112206828247
$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.3102563333333333
$ws.Range("H15").Value = 0.930769
$ws.Range("I15").Value = 0.0005207651961489113
$ws.Range("J15").Value = 0.0005230593158758034
$ws.Range("O15").Value = 0.4367350033185243
$ws.Range("P15").Value = 0.4872798261741802
$ws.Range("Q15").Value = 4.705335865011445
$ws.Range("R15").Value = 42.34802278510301
$ws.Range("S15").Value = 0.0002274363896682667
$ws.Range("T15").Value = 0.0002548762525187471
$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.3102563333333333
$ws.Range("H16").Value = 0.930769
$ws.Range("I16").Value = 0.0005207651961489113
$ws.Range("J16").Value = 0.0005230593158758034
$ws.Range("M16").Value = 10.806162
$ws.Range("N16").Value = 21.612324
$ws.Range("O16").Value = 0.311185607164384
$ws.Range("P16").Value = 0.2314667817857137
$ws.Range("Q16").Value = 3.352680199526
$ws.Range("R16").Value = 20.116081197156
$ws.Range("S16").Value = 0.0001620546337536785
$ws.Range("T16").Value = 0.0001210708565288093

$ws.Range("A17:T21").Delete()
